# Correção nos dados e início da análise PNAD 2009
#
# The sheet had two purely-label "section header" rows (row 5
# "situação do domicílio" and, after the first deletion, what becomes row 7
# "grandes regiões e unidades da federação") that carried no data values of
# their own - the real numeric rows below them were off-by-one (and then
# off-by-two) relative to their intended labels. Deleting those two header
# rows shifts all the data (and the remaining labels) up so each label lines
# up with the correct row of figures, and drops the two now out-of-range
# rows (38/39 -> 39/40 previously "goiás"/"distrito federal" tail) off the
# bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 header: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# Remove the "situação do domicílio" section-header row (row 5).
$ws.Rows(5).Delete()

# After the row-5 deletion, "grandes regiões e unidades da federação"
# (originally row 8) is now at row 7. Remove it too.
$ws.Rows(7).Delete()
